$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 46 entirely; Excel automatically shifts all subsequent rows up by one.
$ws.Rows.Item(46).Delete()

# Append the new work entry as the new last row (225), reusing the same
# ISWC / title / flags that used to belong to the removed row 46.
$ws.Range("A225").Value = 41590407
$ws.Range("B225").Value = "T0400752422"
$ws.Range("C225").Value = "COLO DE ALGODAO"
$ws.Range("D225").Value = "Y"
$ws.Range("E225").Value = "Y"

# Update the view state to reflect where the edit took place.
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("A46:XFD46").Select()
